$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: in the "This service should register a ServiceListener on the
# following: ..." paragraph, the highlighted "BondStreamingService" run goes
# from a black highlight to a red highlight, and the comma that used to start
# the following ", and " run is folded into the highlighted run instead
# (so it reads "...BondStreamingService, and BondInquiryService..." with the
# comma now highlighted red along with the service name).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "BondStreamingService"
$rng.Find.Highlight = $true
$rng.Find.MatchWildcards = $false
$rng.Find.Forward = $true
$rng.Find.Wrap = 0

if ($rng.Find.Execute()) {
    # Recolor the highlight from black to red.
    $rng.Font.HighlightColorIndex = 6

    # Append the comma to this (now red) run ...
    $rng.InsertAfter(",")

    # ... and strip the now-duplicate leading comma off the text that
    # immediately follows ("', and '" -> "' and '").
    $afterRng = $d.Range($rng.End, $rng.End + 1)
    if ($afterRng.Text -eq ",") {
        $afterRng.Delete()
    }
}

# ---------------------------------------------------------------------------
# Change 2: the "Normal" style's font color switches from the automatic
# color to an explicit near-black RGB value (00000A).
# ---------------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.Color = 655360
